$wb = $excel.ActiveWorkbook

# --- Rename existing sheet and add the new "K-Nearest Neighbor" sheet ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "Naive Bayes"

$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "K-Nearest Neighbor"

# --- Populate the new sheet with the n-dimensional Euclidean example data ---
# Fill order chosen so shared-strings end up interleaved the same way as the
# authored workbook (header, then first data row, then remaining headers, ...).

$ws2.Range("A1").Value = "f1"

$ws2.Range("A2").Value = 0.35
$ws2.Range("B2").Value = 0.91
$ws2.Range("C2").Value = 0.86
$ws2.Range("D2").Value = 0.42
$ws2.Range("E2").Value = 0.71
$ws2.Range("F2").Value = "London"

$ws2.Range("B1").Value = "f2"
$ws2.Range("C1").Value = "f3"
$ws2.Range("D1").Value = "f4"
$ws2.Range("F1").Value = "Class"

$ws2.Range("A3").Value = 0.21
$ws2.Range("B3").Value = 0.12
$ws2.Range("C3").Value = 0.76
$ws2.Range("D3").Value = 0.22
$ws2.Range("E3").Value = 0.92
$ws2.Range("F3").Value = "Leeds"

$ws2.Range("A4").Value = 0.41
$ws2.Range("B4").Value = 0.58
$ws2.Range("C4").Value = 0.73
$ws2.Range("D4").Value = 0.21
$ws2.Range("E4").Value = 0.09
$ws2.Range("F4").Value = "Liverpool"

$ws2.Range("E1").Value = "f5"

$ws2.Range("A5").Value = 0.71
$ws2.Range("B5").Value = 0.34
$ws2.Range("C5").Value = 0.55
$ws2.Range("D5").Value = 0.19
$ws2.Range("E5").Value = 0.8
$ws2.Range("F5").Value = "London"

$ws2.Range("A6").Value = 0.79
$ws2.Range("B6").Value = 0.45
$ws2.Range("C6").Value = 0.79
$ws2.Range("D6").Value = 0.21
$ws2.Range("E6").Value = 0.44
$ws2.Range("F6").Value = "Liverpool"

$ws2.Range("A7").Value = 0.61
$ws2.Range("B7").Value = 0.37
$ws2.Range("C7").Value = 0.34
$ws2.Range("D7").Value = 0.81
$ws2.Range("E7").Value = 0.42
$ws2.Range("F7").Value = "Leeds"

$ws2.Range("A8").Value = 0.78
$ws2.Range("B8").Value = 0.12
$ws2.Range("C8").Value = 0.31
$ws2.Range("D8").Value = 0.83
$ws2.Range("E8").Value = 0.87
$ws2.Range("F8").Value = "London"

$ws2.Range("A9").Value = 0.52
$ws2.Range("B9").Value = 0.23
$ws2.Range("C9").Value = 0.73
$ws2.Range("D9").Value = 0.45
$ws2.Range("E9").Value = 0.78
$ws2.Range("F9").Value = "Liverpool"

$ws2.Range("A10").Value = 0.53
$ws2.Range("B10").Value = 0.17
$ws2.Range("C10").Value = 0.63
$ws2.Range("D10").Value = 0.29
$ws2.Range("E10").Value = 0.72
$ws2.Range("F10").Value = "Leeds"

$ws2.Range("A12").Value = 0.65
$ws2.Range("B12").Value = 0.78
$ws2.Range("C12").Value = 0.21
$ws2.Range("D12").Value = 0.29
$ws2.Range("E12").Value = 0.58

# --- View state: new sheet becomes the active/selected one, A12 selected ---
$ws2.Range("A12").Select()
$ws2.Activate()

# --- Workbook window position tweak captured in the authored diff ---
$wb.Windows.Item(1).Left = 21012
$wb.Windows.Item(1).Top = 2496
